$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "Elmar Qarayev"
$ws.Range("B16").Value = "elmarqarayev69@gmail.com"
$ws.Range("C16").Value = 222
$ws.Range("D16").Value = "Pending"
